# Integração com o RocketChat
# The "Troca de modem" ticket in row 2 is replaced: its row is removed
# (shifting the rest of the table up, which carries the little style-only
# marker cells in columns F/G up with it), the header due-date is
# refreshed, and two fresh backlog rows are appended - including the new
# RocketChat ticket with a wrapped-text id cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# drop the old row 2 entirely; rows 3-9 (and their F3/G9 style markers)
# shift up to become rows 2-8
$ws.Rows.Item(2).Delete()

# --- Row 1: due-date refreshed ---
$ws.Range("D1").Value = 43791

# --- Row 2 (used to be row 3): new date for the same ticket/person ---
$ws.Range("D2").Value = 43793.5833333333

# --- Row 3: previously-blank spacer row, now a new backlog entry ---
$ws.Range("A3").Value = "5d8b9232da09856b00edc270"
$ws.Range("B3").Value = "5d825b88b6c09ea9e6d9c732"
$ws.Range("C3").Value = "Troca de modem"
$ws.Range("D3").Value = 43791.5833333333

# --- Row 4: brand-new backlog entry (RocketChat ticket), wrapped id cell ---
$ws.Range("A4").Value = "5d8b9237da09856b00edc276 "
$ws.Range("B4").Value = "5d825b88b6c09ea9e6d9c732"
$ws.Range("C4").Value = "Troca de modem"
$ws.Range("D4").Value = 43791.5833333333
$ws.Range("A4").WrapText = $true

# column D got noticeably narrower once the new rows/ids settled in
$ws.Columns.Item(4).ColumnWidth = 11.2966666667

# keep the sheet's "used range" extending to the very last row (touching a
# cell there is what nudges the used-range/dimension down that far), and
# park the cursor where the author left it
$ws.Range("G1048576").Font.Underline = $true
$ws.Range("G1048576").Font.Underline = $false
$ws.Rows.Item(1048576).RowHeight = 12.8
$ws.Range("C10").Select()
